# Adds two new paragraphs after the existing "=> to Generate Icons" paragraph:
#   "Npm Install react-router"
#   "Npm install react-router-dom"
# and relocates the hidden "_GoBack" bookmark from the first paragraph to the
# end of the (new) last paragraph, matching the target diff.

$d = $word.ActiveDocument

# --- Step 1: insert a new (empty) paragraph right after paragraph 1, then
# fill it with the first new line of text. Using InsertParagraphAfter on a
# range collapsed just before the paragraph mark keeps the original runs in
# paragraph 1 untouched (no run-merging), which matches the diff exactly.
$p1Rng = $d.Paragraphs(1).Range
[void]$p1Rng.MoveEnd(1, -1)
$p1Rng.Collapse(0)
$p1Rng.InsertParagraphAfter()

$p2Rng = $d.Paragraphs(2).Range
$p2Rng.InsertAfter("Npm Install react-router")

# --- Step 2: insert another new paragraph after paragraph 2, and fill it
# with the second new line of text.
$p2RngEnd = $d.Paragraphs(2).Range
[void]$p2RngEnd.MoveEnd(1, -1)
$p2RngEnd.Collapse(0)
$p2RngEnd.InsertParagraphAfter()

$p3Rng = $d.Paragraphs(3).Range
$p3Rng.InsertAfter("Npm install react-router-dom")

# --- Step 3: move the "_GoBack" bookmark from paragraph 1 to the end of the
# new last paragraph (paragraph 3), right after its text, matching the diff.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Work around a COM-bridge quirk where a bookmark collapsed exactly at the
# true end of the document's story gets mis-anchored: insert a throwaway
# buffer character after paragraph 3's text, add the bookmark just before
# that buffer character (so it is no longer sitting at the absolute story
# end), then delete the buffer character again without touching the new
# true end of story (delete range ends one short of it).
$tailRng = $d.Paragraphs(3).Range
[void]$tailRng.MoveEnd(1, -1)
$tailRng.Collapse(0)
$bufferStart = $tailRng.Start
$tailRng.InsertAfter("X")

$bmRng = $d.Range($bufferStart, $bufferStart)
$d.Bookmarks.Add("_GoBack", $bmRng)

$bufferRng = $d.Range($bufferStart, $bufferStart + 1)
$bufferRng.Delete()
